# Auto-generated Excel COM-interop script applying the cryptos.xlsx data refresh
# described by the commit diff (Tue Jul 18 03:54:37 UTC 2023, GitHub Actions).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.121.17"

$ws.Range("D3").Value = "1.910.85"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7386"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.39%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "245.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.59%  "

$ws.Range("E7").Value = "  +0.07%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3096"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.65%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.59"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.22%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06989"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08065"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7722"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.06%  "

$ws.Range("D13").Value = "1.903.89"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.358"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.66%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.18"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.29%  "

$ws.Range("E16").Value = "  -0.55%  "

$ws.Range("D17").Value = "30.120.16"
$ws.Range("E17").Value = "  -0.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.037"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007862"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.32%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.47"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.95%  "

$ws.Range("D21").Value = "2.174.22"
$ws.Range("E21").Value = "  -0.60%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("E23").Value = "  +0.08%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.113"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.99%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.428"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.95%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1286"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.89%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.059"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -6.93%  "

$ws.Range("E30").Value = "  +2.39%  "

$ws.Range("E31").Value = "  -1.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.350"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.103"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.90%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.314"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.05160"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.32%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7516"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.30%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.733"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.21%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01954"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.19%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.799"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.354"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4525"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.78%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "74.60"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.97%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.992"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.04%  "

$ws.Range("E44").Value = "  +0.22%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8402"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.795"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.77%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.963"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.22%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "101.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.57%  "

$ws.Range("D49").Value = "2.073.06"
$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("E50").Value = "  -2.18%  "

$ws.Range("E51").Value = "  -2.66%  "
